$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-06 Tuesday" "2025-05-07 Wednesday"

Replace-Text "215÷5=" "638÷9="
Replace-Text "632÷5=" "282÷6="
Replace-Text "889÷7=" "615÷5="
Replace-Text "376÷8=" "123÷4="
Replace-Text "120÷5=" "537÷8="

Replace-Text "899÷3=" "765÷6="
Replace-Text "440÷3=" "723÷4="
Replace-Text "844÷8=" "821÷5="
Replace-Text "389÷4=" "150÷5="
Replace-Text "361÷6=" "713÷3="

Replace-Text "255÷6=" "564÷2="
Replace-Text "620÷3=" "943÷9="
Replace-Text "497÷9=" "486÷4="
Replace-Text "628÷7=" "710÷9="
Replace-Text "236÷2=" "795÷3="

Replace-Text "865÷4=" "689÷2="
Replace-Text "823÷3=" "791÷2="
Replace-Text "166÷6=" "632÷7="
Replace-Text "525÷6=" "897÷4="
Replace-Text "676÷3=" "839÷3="

Replace-Text "807÷6=" "922÷7="
Replace-Text "997÷6=" "791÷5="
Replace-Text "804÷4=" "143÷5="
Replace-Text "262÷6=" "335÷8="
Replace-Text "211÷7=" "726÷2="
